$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original wide-format data (A1:E7) before rewriting the sheet.
$products = @("Abacaxi*", "Batata-doce", "Cana-de-açúcar", "Mandioca", "Melancia", "Tomate")

$catNames = @("Quantidade 2023/2010", "Valor 2023/2010", "Quantidade 2023/2022", "Valor 2023/2022")

# Columns B,C,D,E hold the values for each category, one per product (rows 2..7)
$values = @{}
for ($col = 2; $col -le 5; $col++) {
    $catIndex = $col - 2
    $catVals = @()
    for ($row = 2; $row -le 7; $row++) {
        $catVals += $ws.Cells.Item($row, $col).Value()
    }
    $values[$catIndex] = $catVals
}

# Clear old contents in the used range
$ws.Range("A1:E7").Clear() | Out-Null

# Write new header
$ws.Cells.Item(1, 1).Value = "Produto"
$ws.Cells.Item(1, 2).Value = "Categoria"
$ws.Cells.Item(1, 3).Value = "Valor"

# Write long-format data: for each category block (6 rows), then next category
$targetRow = 2
for ($catIndex = 0; $catIndex -lt 4; $catIndex++) {
    for ($p = 0; $p -lt 6; $p++) {
        $ws.Cells.Item($targetRow, 1).Value = $products[$p]
        $ws.Cells.Item($targetRow, 2).Value = $catNames[$catIndex]
        $ws.Cells.Item($targetRow, 3).Value = $values[$catIndex][$p]
        $targetRow++
    }
}
